# dangtq them don 3d
# Add a new "3d" order row to the "danh sách đơn" (order list) sheet and
# make that sheet the active tab, matching the author's edit.

$wb = $excel.ActiveWorkbook

$wsOrders = $wb.Worksheets.Item("danh sách đơn")

# New row 8: 3d / OK / 2018-06-13 / path note
$wsOrders.Cells.Item(8, 1).Value = "3d"
$wsOrders.Cells.Item(8, 2).Value = "OK"
$wsOrders.Cells.Item(8, 3).Value = "2018-06-13"
$wsOrders.Cells.Item(8, 4).Value = "\Business\Application forms\3d_Request for appeal"

# This sheet becomes the active/selected tab.
$wsOrders.Activate()
